$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update "employer number" column (A) and "member no" column (B) for rows 2-7
for ($row = 2; $row -le 7; $row++) {
    $ws.Cells.Item($row, 1).Value = 10102345
    $ws.Cells.Item($row, 2).Value = 12345680
}

# Update the active selection to G11
$ws.Range("G11").Select()
